$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A / B data rewrite -------------------------------------------
$rows = @(
    @("COMMENTS", "OFFENSIVE (Y or N)"),
    @("bobo tarantado", "Y"),
    @("gago punyeta", "Y"),
    @("puta fuck gago", "Y"),
    @("maganda", "N"),
    @("tae bobo", "Y"),
    @("pogi", "N"),
    @("matalino mabait", "N"),
    @("masipag magalang atin", "N"),
    @("Yan ung sunod na magdadala Ng pandemia sa", "N"),
    @("Naawa ako sa bata at sa magulang niya", "N"),
    @("Nakakatakot naman ang ginawa ni tatay at dumugo pa ang kamay nya fuck", "Y"),
    @("bless kuya", $null)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    if ($null -eq $rows[$i][1]) {
        $ws.Cells.Item($r, 2).ClearContents()
    } else {
        $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    }
}

# --- Column width (author typed "100" in Format Cells > Column Width;
#     Excel's MDW-7 pixel-grid formula stores that as 100.7109375) ----------
$ws.Range("A:B").ColumnWidth = 100

# --- Selection ---------------------------------------------------------------
$ws.Range("A13").Select()
